# Cotações atualizadas - 2025-11-15
# Append a new row (row 72) with the quote values for 2025-11-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 72

# Column A: date serial, formatted like the existing date column (copy format from A71)
$ws.Cells.Item($row, 1).Value = 45976
$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat

# Columns B:E: textual quote values (kept as text, matching the rest of the column)
$ws.Cells.Item($row, 2).Value = "22,2279"
$ws.Cells.Item($row, 3).Value = "15,8858"
$ws.Cells.Item($row, 4).Value = "15,6322"
$ws.Cells.Item($row, 5).Value = "15,6322"
